$wb = $excel.ActiveWorkbook

# Sheet ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2076.4119
$ws.Range("I137").Value = 1864.0714
$ws.Range("K137").Value = 5592.2142
$ws.Range("M137").Value = -3042.2142

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 50005284
$ws.Range("J138").Value = 4758.8887
$ws.Range("L138").Value = 14276.6661
$ws.Range("N138").Value = -24556.6661

# Sheet ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2242.4666
$ws.Range("I141").Value = 1571.5454
$ws.Range("J141").Value = 4087.5
$ws.Range("K141").Value = 4714.6362
$ws.Range("L141").Value = 12262.5
$ws.Range("M141").Value = 465.3638000000001
$ws.Range("N141").Value = -22622.5

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4812.8945
$ws.Range("I32").Value = 3736.46
$ws.Range("J32").Value = 12501.714
$ws.Range("K32").Value = 3736.46
$ws.Range("L32").Value = 12501.714
$ws.Range("M32").Value = -3449.46
$ws.Range("N32").Value = -13075.714

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3507.1667
$ws.Range("J45").Value = 3861.1428
$ws.Range("L45").Value = 3861.1428
$ws.Range("N45").Value = -4615.1428

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 71431700
$ws.Range("I74").Value = 142859400
$ws.Range("J74").Value = 3994
$ws.Range("K74").Value = 142859400
$ws.Range("L74").Value = 3994
$ws.Range("M74").Value = -142858526
$ws.Range("N74").Value = -5742

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 71431700
$ws.Range("I77").Value = 142859400
$ws.Range("J77").Value = 3994
$ws.Range("K77").Value = 714297000
$ws.Range("L77").Value = 19970
$ws.Range("M77").Value = -714292632
$ws.Range("N77").Value = -28706

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2697.3215
$ws.Range("I122").Value = 2114.8635
$ws.Range("K122").Value = 6344.5905
$ws.Range("M122").Value = -3894.5905

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 18533.709
$ws.Range("I132").Value = 2217.0435
$ws.Range("K132").Value = 6651.130500000001
$ws.Range("M132").Value = -4121.130500000001

# Sheet BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 779.1111
$ws.Range("J107").Value = 729.6667
$ws.Range("L107").Value = 729.6667
$ws.Range("N107").Value = -4569.6667

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3646.9656
$ws.Range("I134").Value = 4035.92
$ws.Range("K134").Value = 12107.76
$ws.Range("M134").Value = -9572.76

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20973.96
$ws.Range("I58").Value = 1726.8182
$ws.Range("J58").Value = 35088.535
$ws.Range("K58").Value = 1726.8182
$ws.Range("L58").Value = 35088.535
$ws.Range("M58").Value = -1523.8182
$ws.Range("N58").Value = -35494.535

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2220
$ws.Range("I122").Value = 2366.6667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7100.000100000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4650.000100000001
$ws.Range("N122").Value = -10900

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3866.2
$ws.Range("I132").Value = 2919.077
$ws.Range("J132").Value = 5625.143
$ws.Range("K132").Value = 8757.231
$ws.Range("L132").Value = 16875.429
$ws.Range("M132").Value = -6227.231
$ws.Range("N132").Value = -21935.429

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1022.89655
$ws.Range("I134").Value = 910.26086
$ws.Range("J134").Value = 1454.6666
$ws.Range("K134").Value = 2730.78258
$ws.Range("L134").Value = 4363.9998
$ws.Range("M134").Value = -195.7825800000001
$ws.Range("N134").Value = -9433.9998

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 20973.96
$ws.Range("I136").Value = 1726.8182
$ws.Range("J136").Value = 35088.535
$ws.Range("K136").Value = 5180.4546
$ws.Range("L136").Value = 105265.605
$ws.Range("M136").Value = -2630.4546
$ws.Range("N136").Value = -110365.605

# Sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1197.05
$ws.Range("I5").Value = 820.1667
$ws.Range("J5").Value = 1762.375
$ws.Range("K5").Value = 2460.5001
$ws.Range("L5").Value = 5287.125
$ws.Range("M5").Value = -2348.5001
$ws.Range("N5").Value = -5511.125

# Sheet CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3470.6667
$ws.Range("J58").Value = 3964.8
$ws.Range("L58").Value = 11894.4
$ws.Range("N58").Value = -12150.4

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 699.66
$ws.Range("J131").Value = 719.8495
$ws.Range("L131").Value = 2159.5485
$ws.Range("N131").Value = -12239.5485

# Sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1197.05
$ws.Range("I135").Value = 820.1667
$ws.Range("J135").Value = 1762.375
$ws.Range("K135").Value = 7381.5003
$ws.Range("L135").Value = 15861.375
$ws.Range("M135").Value = -4846.5003
$ws.Range("N135").Value = -20931.375

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1154.9667
$ws.Range("I97").Value = 1343.5625
$ws.Range("J97").Value = 939.4286
$ws.Range("K97").Value = 1343.5625
$ws.Range("L97").Value = 939.4286
$ws.Range("M97").Value = -847.5625
$ws.Range("N97").Value = -1931.4286

# Sheet GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5495.2173
$ws.Range("J126").Value = 6789
$ws.Range("L126").Value = 20367
$ws.Range("N126").Value = -25307

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27538.1
$ws.Range("I132").Value = 2415.6
$ws.Range("K132").Value = 7246.799999999999
$ws.Range("M132").Value = -4716.799999999999

# Sheet LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10089.333
$ws.Range("J7").Value = 16676
$ws.Range("L7").Value = 16676
$ws.Range("N7").Value = -16900

# Sheet LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 10089.333
$ws.Range("J126").Value = 16676
$ws.Range("L126").Value = 50028
$ws.Range("N126").Value = -54968

# Sheet LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1613.7391
$ws.Range("I136").Value = 1423.1111
$ws.Range("K136").Value = 4269.3333
$ws.Range("M136").Value = -1719.3333

# Sheet WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1662.65
$ws.Range("I132").Value = 879.5833
$ws.Range("J132").Value = 2837.25
$ws.Range("K132").Value = 2638.7499
$ws.Range("L132").Value = 8511.75
$ws.Range("M132").Value = -108.7498999999998
$ws.Range("N132").Value = -13571.75

# Sheet WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 22941648
$ws.Range("I136").Value = 30361662
$ws.Range("J136").Value = 7064
$ws.Range("K136").Value = 91084986
$ws.Range("L136").Value = 21192
$ws.Range("M136").Value = -91082436
$ws.Range("N136").Value = -26292
